# "take out lignite from ets"
# Zero-out the "lignite" shareweight row on the ETS sheet (row 13, years
# 2020-2050 in columns B:AF), and leave the workbook with the ETS sheet
# active/selected (cell B14), matching the author's final UI state.

$wb = $excel.ActiveWorkbook

$ets = $wb.Worksheets.Item("ETS")

# Lignite shareweights (row 13, columns B:AF) -> 0 for every year.
$ets.Range("B13:AF13").Value = 0

# Switch focus to the ETS sheet and leave the cursor on B14, as in the
# committed workbook.
$ets.Activate()
$ets.Range("B14").Select()
